$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 58
$ws.Range("I2").Value = 161
$ws.Range("J2").Value = 737
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 198
$ws.Range("M2").Value = 18
$ws.Range("N2").Value = 135
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 13
$ws.Range("S2").Value = 70
$ws.Range("T2").Value = 115
$ws.Range("U2").Value = 8
$ws.Range("V2").Value = 1152
$ws.Range("X2").Value = 1110
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 9
